# The deck's theme colours were accidentally swapped between the two
# embedded theme parts: the slides (via the single Slide Master) were
# rendering the "Integral" palette while the intended "Office Theme"
# palette lived, unused, in the other theme part. Restore the deck's
# live/active theme (reachable from every slide through the shared
# Slide Master) to the "Office Theme" colour palette.
#
# PowerPoint's object model exposes the 12 theme colour slots (dk1,
# lt1, dk2, lt2, accent1-6, hlink, folHlink) as a ThemeColorScheme on
# a Slide/SlideRange; writing .RGB there rewrites the colours of the
# theme part actually bound to the presentation's Slide Master.
# .RGB uses the Windows COLORREF (0xBBGGRR) ordering.

$p = $ppt.ActivePresentation
$tcs = $p.Slides.Item(1).ThemeColorScheme

$tcs.Item(1).RGB  = 0x000000   # dk1      -> 000000
$tcs.Item(2).RGB  = 0xFFFFFF   # lt1      -> FFFFFF
$tcs.Item(3).RGB  = 0x6A5444   # dk2      -> 44546A
$tcs.Item(4).RGB  = 0xE6E6E7   # lt2      -> E7E6E6
$tcs.Item(5).RGB  = 0xD59B5B   # accent1  -> 5B9BD5
$tcs.Item(6).RGB  = 0x317DED   # accent2  -> ED7D31
$tcs.Item(7).RGB  = 0xA5A5A5   # accent3  -> A5A5A5
$tcs.Item(8).RGB  = 0x00C0FF   # accent4  -> FFC000
$tcs.Item(9).RGB  = 0xC47244   # accent5  -> 4472C4
$tcs.Item(10).RGB = 0x47AD70   # accent6  -> 70AD47
$tcs.Item(11).RGB = 0xC16305   # hlink    -> 0563C1
$tcs.Item(12).RGB = 0x724F95   # folHlink -> 954F72
